$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve an empty-but-present cell template from AF50 before it gets cleared,
# then use it to stamp AF48 as an empty (but present) cell, matching the source XML.
$ws.Range("AF50").Copy($ws.Range("AF48"))

# Row 44
$ws.Range("A44").Value = 111684838
$ws.Range("B44").Value = 90682
$ws.Range("D44").Value = 'NT'
$ws.Range("E44").Value = 2059
$ws.Range("F44").Value = 'Skrovlig taggsvamp'
$ws.Range("G44").Value = 'Hydnellum scabrosum'
$ws.Range("H44").Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("Q44").Value = 538321
$ws.Range("R44").Value = 7007201
$ws.Range("Z44").ClearContents()
$ws.Range("AB44").ClearContents()
$ws.Range("AC44").ClearContents()
$ws.Range("AF44").ClearContents()

# Row 45
$ws.Range("A45").Value = 111684866
$ws.Range("B45").Value = 90678
$ws.Range("D45").Value = 'LC'
$ws.Range("E45").Value = 4366
$ws.Range("F45").Value = 'Skarp dropptaggsvamp'
$ws.Range("G45").Value = 'Hydnellum peckii'
$ws.Range("H45").Value = 'Banker'
$ws.Range("Q45").Value = 538378
$ws.Range("R45").Value = 7007068
$ws.Range("Z45").ClearContents()
$ws.Range("AB45").ClearContents()
$ws.Range("AC45").ClearContents()
$ws.Range("AF45").ClearContents()

# Row 46
$ws.Range("A46").Value = 111684869
$ws.Range("B46").Value = 78578
$ws.Range("D46").Value = 'NT'
$ws.Range("E46").Value = 6458
$ws.Range("F46").Value = 'Lunglav'
$ws.Range("G46").Value = 'Lobaria pulmonaria'
$ws.Range("H46").Value = '(L.) Hoffm.'
$ws.Range("Q46").Value = 538403
$ws.Range("R46").Value = 7007022
$ws.Range("Z46").ClearContents()
$ws.Range("AB46").ClearContents()
$ws.Range("AC46").ClearContents()
$ws.Range("AF46").ClearContents()

# Row 47
$ws.Range("A47").Value = 111684854
$ws.Range("B47").Value = 96346
$ws.Range("D47").Value = 'NT'
$ws.Range("E47").Value = 620
$ws.Range("F47").Value = 'Skogsfru'
$ws.Range("G47").Value = 'Epipogium aphyllum'
$ws.Range("H47").Value = 'Sw.'
$ws.Range("Q47").Value = 538430
$ws.Range("R47").Value = 7007062
$ws.Range("Z47").ClearContents()
$ws.Range("AB47").ClearContents()
$ws.Range("AC47").ClearContents()
$ws.Range("AF47").ClearContents()

# Row 48
$ws.Range("A48").Value = 111684835
$ws.Range("B48").Value = 85266
$ws.Range("D48").Value = 'LC'
$ws.Range("E48").Value = 249228
$ws.Range("F48").Value = 'Barrfagerspindling'
$ws.Range("G48").Value = 'Cortinarius piceae'
$ws.Range("H48").Value = 'Frøslev, T.S.Jeppesen & Brandrud'
$ws.Range("Q48").Value = 538475
$ws.Range("R48").Value = 7007186
$ws.Range("Z48").ClearContents()
$ws.Range("AB48").ClearContents()
$ws.Range("AC48").Value = 'Både gran och tall. Gul hatt med blek kant. Mörka velumfläckar i mitten. Gulaktig på foten. Ingen reaktion KOH på hatten. På foten ingen eller mörkbrun.'

# Row 49
$ws.Range("A49").Value = 111684843
$ws.Range("B49").Value = 90332
$ws.Range("D49").Value = 'LC'
$ws.Range("E49").Value = 4769
$ws.Range("F49").Value = 'Svavelriska'
$ws.Range("G49").Value = 'Lactarius scrobiculatus'
$ws.Range("H49").Value = '(Scop.:Fr.) Fr.'
$ws.Range("Q49").Value = 538471
$ws.Range("R49").Value = 7007183
$ws.Range("Z49").ClearContents()
$ws.Range("AB49").ClearContents()
$ws.Range("AC49").ClearContents()
$ws.Range("AF49").ClearContents()

# Row 50
$ws.Range("A50").Value = 111684853
$ws.Range("B50").Value = 96253
$ws.Range("D50").Value = 'LC'
$ws.Range("E50").Value = 504
$ws.Range("F50").Value = 'Guckusko'
$ws.Range("G50").Value = 'Cypripedium calceolus'
$ws.Range("H50").Value = 'L.'
$ws.Range("Q50").Value = 538430
$ws.Range("R50").Value = 7007062
$ws.Range("Z50").ClearContents()
$ws.Range("AB50").ClearContents()
$ws.Range("AC50").ClearContents()
$ws.Range("AF50").ClearContents()

# Row 51
$ws.Range("A51").Value = 111684868
$ws.Range("B51").Value = 90666
$ws.Range("D51").Value = 'LC'
$ws.Range("E51").Value = 4364
$ws.Range("F51").Value = 'Dropptaggsvamp'
$ws.Range("G51").Value = 'Hydnellum ferrugineum'
$ws.Range("H51").Value = '(Fr.:Fr.) P. Karst.'
$ws.Range("Q51").Value = 538403
$ws.Range("R51").Value = 7007022
$ws.Range("Z51").ClearContents()
$ws.Range("AB51").ClearContents()
$ws.Range("AC51").ClearContents()
$ws.Range("AF51").ClearContents()

# Row 52
$ws.Range("A52").Value = 111684865
$ws.Range("B52").Value = 78578
$ws.Range("D52").Value = 'NT'
$ws.Range("E52").Value = 6458
$ws.Range("F52").Value = 'Lunglav'
$ws.Range("G52").Value = 'Lobaria pulmonaria'
$ws.Range("H52").Value = '(L.) Hoffm.'
$ws.Range("Q52").Value = 538378
$ws.Range("R52").Value = 7007068
$ws.Range("Z52").ClearContents()
$ws.Range("AB52").ClearContents()
$ws.Range("AC52").ClearContents()
$ws.Range("AF52").ClearContents()

# Row 53
$ws.Range("A53").Value = 111684880
$ws.Range("B53").Value = 86223
$ws.Range("D53").Value = 'NT'
$ws.Range("E53").Value = 4412
$ws.Range("F53").Value = 'Äggvaxskivling'
$ws.Range("G53").Value = 'Hygrophorus karstenii'
$ws.Range("H53").Value = 'Sacc. & Cub.'
$ws.Range("Q53").Value = 538506
$ws.Range("R53").Value = 7007072
$ws.Range("Z53").ClearContents()
$ws.Range("AB53").ClearContents()
$ws.Range("AC53").ClearContents()
$ws.Range("AF53").ClearContents()

# Row 54
$ws.Range("A54").Value = 111684878
$ws.Range("B54").Value = 86223
$ws.Range("D54").Value = 'NT'
$ws.Range("E54").Value = 4412
$ws.Range("F54").Value = 'Äggvaxskivling'
$ws.Range("G54").Value = 'Hygrophorus karstenii'
$ws.Range("H54").Value = 'Sacc. & Cub.'
$ws.Range("Q54").Value = 538418
$ws.Range("R54").Value = 7007018
$ws.Range("Z54").ClearContents()
$ws.Range("AB54").ClearContents()
$ws.Range("AC54").ClearContents()
$ws.Range("AF54").ClearContents()

# Row 55
$ws.Range("A55").Value = 111684892
$ws.Range("B55").Value = 96253
$ws.Range("D55").Value = 'LC'
$ws.Range("E55").Value = 504
$ws.Range("F55").Value = 'Guckusko'
$ws.Range("G55").Value = 'Cypripedium calceolus'
$ws.Range("H55").Value = 'L.'
$ws.Range("Q55").Value = 538447
$ws.Range("R55").Value = 7007049
$ws.Range("Z55").ClearContents()
$ws.Range("AB55").ClearContents()
$ws.Range("AC55").ClearContents()
$ws.Range("AF55").ClearContents()

# Row 56
$ws.Range("A56").Value = 111684896
$ws.Range("B56").Value = 90710
$ws.Range("D56").Value = 'NT'
$ws.Range("E56").Value = 5449
$ws.Range("F56").Value = 'Svart taggsvamp'
$ws.Range("G56").Value = 'Phellodon niger'
$ws.Range("H56").Value = '(Fr.:Fr.) P.Karst.'
$ws.Range("Q56").Value = 538349
$ws.Range("R56").Value = 7007238
$ws.Range("Z56").ClearContents()
$ws.Range("AB56").ClearContents()
$ws.Range("AC56").ClearContents()
$ws.Range("AF56").ClearContents()

# Row 57
$ws.Range("A57").Value = 111684888
$ws.Range("B57").Value = 90671
$ws.Range("D57").Value = 'VU'
$ws.Range("E57").Value = 6003298
$ws.Range("F57").Value = 'Ruttaggsvamp'
$ws.Range("G57").Value = 'Hydnellum illudens'
$ws.Range("H57").Value = '(Maas Geest.) Nitare'
$ws.Range("Q57").Value = 538472
$ws.Range("R57").Value = 7007185
$ws.Range("Z57").ClearContents()
$ws.Range("AB57").ClearContents()
$ws.Range("AC57").Value = 'Kalktallskog med enbuskar och stenbär. Gult kött vid torkning.'
